# Fixed StudyComb for Faceted Filters ICDC
#
# The StatQuery column (C) for each tab row (CasesTab, SamplesTab, FilesTab)
# previously shared one long Cypher query used to compute summary statistics.
# That query is replaced with a new, corrected query.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStatQuery = "MATCH (demo:demographic)`nWHERE demo.breed IN ['Australian Shepherd']`nMATCH (demo:demographic)-->(c:case)-->(s:study)-->(p:program)`nOPTIONAL MATCH (c)<-[*]-(samp:sample)`nOPTIONAL MATCH (c)<-[*]-(f:file)`nRETURN `n`tcount(DISTINCT(f)) as number_of_files, `n`tcount(DISTINCT(samp)) as number_of_sample, `n`tcount(DISTINCT(c)) as number_of_cases, `n`tcount(DISTINCT(s)) as number_of_study"

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Update the saved view state: scroll so row 4 / column A is the top-left
# visible cell, set zoom to 100%, and move the active selection to B4.
$win = $excel.ActiveWindow
$ws.Range("B4").Select()
$win.ScrollRow = 4
$win.ScrollColumn = 1
$win.Zoom = 100
